$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# E1 currently holds the date serial for 2024-03-03; replace it with the
# text label used by the other header cells (B1/C1/D1), matching the
# "dd_mm_yyyy" string style.
$ws.Range("E1").Value = "03_03_2024"

# Move the active selection from E6 to E2.
$ws.Range("E2").Select()
